$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.23349936394030379,
    0.015525888577602842,
    -0.0039999999863447044,
    -0.00799999997501466,
    -0.0029999999850174319,
    -0.0019999999830666582,
    -0.0099999999638895609,
    0.019042946294359098,
    -0.001999999981010081,
    -0.0019999999796347367,
    -0.0029999999770842223,
    -0.0034999999755345179,
    -0.0034999999741396337,
    -0.0079999999630580021,
    -0.00099999997938926555,
    -0.0019999999767597032,
    -0.0019999999763973264,
    -0.0039999999715565338,
    -0.0039999999895043992,
    -0.0039999999863340463,
    -0.047819312204852693,
    -0.0039999999854893886,
    -0.0049999999832239794,
    -0.019999999946339386,
    -0.019999999945681246,
    -0.002499999981134593,
    -0.0024999999800674466,
    -0.0019999999770012877,
    -0.0069999999619563269,
    -0.059999999834727813,
    -0.0069999999592269546,
    -0.009999999951928018,
    -0.0039999999660942365
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
